# This script updates the "F" column ("想去人数" / want-to-go count) values
# on several rows across the four worksheets, matching the values produced
# by a newer scrape of the source data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 8415
$ws1.Range("F9").Value  = 23
$ws1.Range("F11").Value = 7675
$ws1.Range("F22").Value = 966
$ws1.Range("F23").Value = 1506
$ws1.Range("F29").Value = 15
$ws1.Range("F33").Value = 1304
$ws1.Range("F34").Value = 475
$ws1.Range("F37").Value = 236
$ws1.Range("F41").Value = 2511

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value  = 133
$ws2.Range("F9").Value  = 22
$ws2.Range("F45").Value = 52

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value  = 1645
$ws3.Range("F7").Value  = 696
$ws3.Range("F9").Value  = 9494
$ws3.Range("F10").Value = 1821
$ws3.Range("F11").Value = 198

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 8415
$ws4.Range("F6").Value  = 696
$ws4.Range("F7").Value  = 1821
$ws4.Range("F17").Value = 967
$ws4.Range("F18").Value = 1506
$ws4.Range("F19").Value = 133
$ws4.Range("F22").Value = 22
$ws4.Range("F32").Value = 1304
$ws4.Range("F35").Value = 475
$ws4.Range("F39").Value = 236
$ws4.Range("F46").Value = 2511
$ws4.Range("F48").Value = 52
